# adding auto Define Time (PAGI/MALAM)
# Sets column J (rows 4-17) to "PAGI" for each data row on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 4; $row -le 17; $row++) {
    $ws.Cells.Item($row, 10).Value = "PAGI"
}
